$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Btc"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.191714
$ws.Range("H2").Value = 0.575142
$ws.Range("I2").Value = 0.09369188973541917
$ws.Range("J2").Value = 0.09369188973541917
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04178033333333334
$ws.Range("N2").Value = 0.125341
$ws.Range("O2").Value = 0.009630623139527362
$ws.Range("P2").Value = 0.009630623139527362
$ws.Range("Q2").Value = 0.008009874824666668
$ws.Range("R2").Value = 0.07208887342200002
$ws.Range("S2").Value = 0.000902311281271974
$ws.Range("T2").Value = 0.000902311281271974

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Btc"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.191714
$ws.Range("H3").Value = 0.575142
$ws.Range("I3").Value = 0.09369188973541917
$ws.Range("J3").Value = 0.09369188973541917
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.181585
$ws.Range("N3").Value = 0.544755
$ws.Range("O3").Value = 0.04185645645377991
$ws.Range("P3").Value = 0.04185645645377991
$ws.Range("Q3").Value = 0.03481238669
$ws.Range("R3").Value = 0.31331148021
$ws.Range("S3").Value = 0.003921610502782921
$ws.Range("T3").Value = 0.003921610502782921

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Btc"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.191714
$ws.Range("H4").Value = 0.575142
$ws.Range("I4").Value = 0.09369188973541917
$ws.Range("J4").Value = 0.09369188973541917
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4121513333333333
$ws.Range("N4").Value = 1.236454
$ws.Range("O4").Value = 0.09500341072243849
$ws.Range("P4").Value = 0.09500341072243849
$ws.Range("Q4").Value = 0.07901518071866667
$ws.Range("R4").Value = 0.711136626468
$ws.Range("S4").Value = 0.008901049081895446
$ws.Range("T4").Value = 0.008901049081895446

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Btc"
$ws.Range("C5").Value = "Erbb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.191714
$ws.Range("H5").Value = 0.575142
$ws.Range("I5").Value = 0.09369188973541917
$ws.Range("J5").Value = 0.09369188973541917
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.702762666666667
$ws.Range("N5").Value = 11.108288
$ws.Range("O5").Value = 0.8535095096842542
$ws.Range("P5").Value = 0.8535095096842543
$ws.Range("Q5").Value = 0.7098714418773334
$ws.Range("R5").Value = 6.388842976896
$ws.Range("S5").Value = 0.07996691886946883
$ws.Range("T5").Value = 0.07996691886946884

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Btc"
$ws.Range("C6").Value = "Erbb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.854503666666667
$ws.Range("H6").Value = 5.563511
$ws.Range("I6").Value = 0.9063081102645809
$ws.Range("J6").Value = 0.9063081102645809
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04178033333333334
$ws.Range("N6").Value = 0.125341
$ws.Range("O6").Value = 0.009630623139527362
$ws.Range("P6").Value = 0.009630623139527362
$ws.Range("Q6").Value = 0.07748178136122223
$ws.Range("R6").Value = 0.697336032251
$ws.Range("S6").Value = 0.008728311858255388
$ws.Range("T6").Value = 0.008728311858255388

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Btc"
$ws.Range("C7").Value = "Erbb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.854503666666667
$ws.Range("H7").Value = 5.563511
$ws.Range("I7").Value = 0.9063081102645809
$ws.Range("J7").Value = 0.9063081102645809
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.181585
$ws.Range("N7").Value = 0.544755
$ws.Range("O7").Value = 0.04185645645377991
$ws.Range("P7").Value = 0.04185645645377991
$ws.Range("Q7").Value = 0.3367500483116667
$ws.Range("R7").Value = 3.030750434805
$ws.Range("S7").Value = 0.03793484595099699
$ws.Range("T7").Value = 0.03793484595099699

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Btc"
$ws.Range("C8").Value = "Erbb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.854503666666667
$ws.Range("H8").Value = 5.563511
$ws.Range("I8").Value = 0.9063081102645809
$ws.Range("J8").Value = 0.9063081102645809
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4121513333333333
$ws.Range("N8").Value = 1.236454
$ws.Range("O8").Value = 0.09500341072243849
$ws.Range("P8").Value = 0.09500341072243849
$ws.Range("Q8").Value = 0.7643361588882223
$ws.Range("R8").Value = 6.879025429994
$ws.Range("S8").Value = 0.08610236164054305
$ws.Range("T8").Value = 0.08610236164054305

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Btc"
$ws.Range("C9").Value = "Erbb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.854503666666667
$ws.Range("H9").Value = 5.563511
$ws.Range("I9").Value = 0.9063081102645809
$ws.Range("J9").Value = 0.9063081102645809
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.702762666666667
$ws.Range("N9").Value = 11.108288
$ws.Range("O9").Value = 0.8535095096842542
$ws.Range("P9").Value = 0.8535095096842543
$ws.Range("Q9").Value = 6.866786942129778
$ws.Range("R9").Value = 61.801082479168
$ws.Range("S9").Value = 0.7735425908147855
$ws.Range("T9").Value = 0.7735425908147856
